$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 1043.7273  # H80
$ws.Cells.Item(80, 10).Value = 1623.625  # J80
$ws.Cells.Item(80, 12).Value = 4870.875  # L80
$ws.Cells.Item(80, 14).Value = -6866.875  # N80
$ws.Cells.Item(82, 8).Value = 500  # H82
$ws.Cells.Item(82, 9).Value = 500  # I82
$ws.Cells.Item(82, 11).Value = 1500  # K82
$ws.Cells.Item(82, 13).Value = -1094  # M82
$ws.Cells.Item(83, 8).Value = 1043.7273  # H83
$ws.Cells.Item(83, 10).Value = 1623.625  # J83
$ws.Cells.Item(83, 12).Value = 14612.625  # L83
$ws.Cells.Item(83, 14).Value = -24596.625  # N83
$ws.Cells.Item(85, 8).Value = 500  # H85
$ws.Cells.Item(85, 9).Value = 500  # I85
$ws.Cells.Item(85, 11).Value = 1500  # K85
$ws.Cells.Item(85, 13).Value = -96  # M85
$ws.Cells.Item(99, 8).Value = 871.5  # H99
$ws.Cells.Item(99, 9).Value = 899.6667  # I99
$ws.Cells.Item(99, 10).Value = 787  # J99
$ws.Cells.Item(99, 11).Value = 2699.0001  # K99
$ws.Cells.Item(99, 12).Value = 2361  # L99
$ws.Cells.Item(99, 13).Value = -1201.0001  # M99
$ws.Cells.Item(99, 14).Value = -5357  # N99
$ws.Cells.Item(106, 8).Value = 41695284  # H106
$ws.Cells.Item(106, 9).Value = 43504864  # I106
$ws.Cells.Item(106, 11).Value = 43504864  # K106
$ws.Cells.Item(106, 13).Value = -43504233  # M106
$ws.Cells.Item(112, 8).Value = 3454.1667  # H112
$ws.Cells.Item(112, 9).Value = 2983.3333  # I112
$ws.Cells.Item(112, 11).Value = 8949.999899999999  # K112
$ws.Cells.Item(112, 13).Value = -7841.999899999999  # M112
$ws.Cells.Item(125, 8).Value = $null  # H125 (cleared)
$ws.Cells.Item(125, 9).Value = $null  # I125 (cleared)
$ws.Cells.Item(125, 10).Value = $null  # J125 (cleared)
$ws.Cells.Item(125, 11).Value = $null  # K125 (cleared)
$ws.Cells.Item(125, 12).Value = $null  # L125 (cleared)
$ws.Cells.Item(125, 13).Value = $null  # M125 (cleared)
$ws.Cells.Item(126, 8).Value = $null  # H126 (cleared)
$ws.Cells.Item(126, 9).Value = $null  # I126 (cleared)
$ws.Cells.Item(126, 10).Value = $null  # J126 (cleared)
$ws.Cells.Item(126, 11).Value = $null  # K126 (cleared)
$ws.Cells.Item(126, 12).Value = $null  # L126 (cleared)
$ws.Cells.Item(127, 8).Value = $null  # H127 (cleared)
$ws.Cells.Item(127, 9).Value = $null  # I127 (cleared)
$ws.Cells.Item(127, 10).Value = $null  # J127 (cleared)
$ws.Cells.Item(127, 11).Value = $null  # K127 (cleared)
$ws.Cells.Item(127, 12).Value = $null  # L127 (cleared)
$ws.Cells.Item(127, 13).Value = $null  # M127 (cleared)
$ws.Cells.Item(128, 8).Value = $null  # H128 (cleared)
$ws.Cells.Item(128, 9).Value = $null  # I128 (cleared)
$ws.Cells.Item(128, 10).Value = $null  # J128 (cleared)
$ws.Cells.Item(128, 11).Value = $null  # K128 (cleared)
$ws.Cells.Item(128, 12).Value = $null  # L128 (cleared)
$ws.Cells.Item(129, 8).Value = $null  # H129 (cleared)
$ws.Cells.Item(129, 9).Value = $null  # I129 (cleared)
$ws.Cells.Item(129, 10).Value = $null  # J129 (cleared)
$ws.Cells.Item(129, 11).Value = $null  # K129 (cleared)
$ws.Cells.Item(129, 12).Value = $null  # L129 (cleared)
$ws.Cells.Item(129, 13).Value = $null  # M129 (cleared)
$ws.Cells.Item(129, 14).Value = $null  # N129 (cleared)
$ws.Cells.Item(130, 8).Value = $null  # H130 (cleared)
$ws.Cells.Item(130, 9).Value = $null  # I130 (cleared)
$ws.Cells.Item(130, 10).Value = $null  # J130 (cleared)
$ws.Cells.Item(130, 11).Value = $null  # K130 (cleared)
$ws.Cells.Item(130, 12).Value = $null  # L130 (cleared)
$ws.Cells.Item(130, 14).Value = $null  # N130 (cleared)
$ws.Cells.Item(131, 8).Value = $null  # H131 (cleared)
$ws.Cells.Item(131, 9).Value = $null  # I131 (cleared)
$ws.Cells.Item(131, 10).Value = $null  # J131 (cleared)
$ws.Cells.Item(131, 11).Value = $null  # K131 (cleared)
$ws.Cells.Item(131, 12).Value = $null  # L131 (cleared)
$ws.Cells.Item(131, 13).Value = $null  # M131 (cleared)
$ws.Cells.Item(131, 14).Value = $null  # N131 (cleared)
$ws.Cells.Item(132, 8).Value = $null  # H132 (cleared)
$ws.Cells.Item(132, 9).Value = $null  # I132 (cleared)
$ws.Cells.Item(132, 10).Value = $null  # J132 (cleared)
$ws.Cells.Item(132, 11).Value = $null  # K132 (cleared)
$ws.Cells.Item(132, 12).Value = $null  # L132 (cleared)
$ws.Cells.Item(132, 13).Value = $null  # M132 (cleared)
$ws.Cells.Item(132, 14).Value = $null  # N132 (cleared)
$ws.Cells.Item(133, 8).Value = $null  # H133 (cleared)
$ws.Cells.Item(133, 9).Value = $null  # I133 (cleared)
$ws.Cells.Item(133, 10).Value = $null  # J133 (cleared)
$ws.Cells.Item(133, 11).Value = $null  # K133 (cleared)
$ws.Cells.Item(133, 12).Value = $null  # L133 (cleared)
$ws.Cells.Item(134, 8).Value = $null  # H134 (cleared)
$ws.Cells.Item(134, 9).Value = $null  # I134 (cleared)
$ws.Cells.Item(134, 10).Value = $null  # J134 (cleared)
$ws.Cells.Item(134, 11).Value = $null  # K134 (cleared)
$ws.Cells.Item(134, 12).Value = $null  # L134 (cleared)
$ws.Cells.Item(135, 8).Value = $null  # H135 (cleared)
$ws.Cells.Item(135, 9).Value = $null  # I135 (cleared)
$ws.Cells.Item(135, 10).Value = $null  # J135 (cleared)
$ws.Cells.Item(135, 11).Value = $null  # K135 (cleared)
$ws.Cells.Item(135, 12).Value = $null  # L135 (cleared)
$ws.Cells.Item(135, 13).Value = $null  # M135 (cleared)
$ws.Cells.Item(135, 14).Value = $null  # N135 (cleared)
$ws.Cells.Item(136, 8).Value = $null  # H136 (cleared)
$ws.Cells.Item(136, 9).Value = $null  # I136 (cleared)
$ws.Cells.Item(136, 10).Value = $null  # J136 (cleared)
$ws.Cells.Item(136, 11).Value = $null  # K136 (cleared)
$ws.Cells.Item(136, 12).Value = $null  # L136 (cleared)
$ws.Cells.Item(137, 8).Value = $null  # H137 (cleared)
$ws.Cells.Item(137, 9).Value = $null  # I137 (cleared)
$ws.Cells.Item(137, 10).Value = $null  # J137 (cleared)
$ws.Cells.Item(137, 11).Value = $null  # K137 (cleared)
$ws.Cells.Item(137, 12).Value = $null  # L137 (cleared)
$ws.Cells.Item(137, 13).Value = $null  # M137 (cleared)
$ws.Cells.Item(137, 14).Value = $null  # N137 (cleared)
$ws.Cells.Item(138, 8).Value = $null  # H138 (cleared)
$ws.Cells.Item(138, 9).Value = $null  # I138 (cleared)
$ws.Cells.Item(138, 10).Value = $null  # J138 (cleared)
$ws.Cells.Item(138, 11).Value = $null  # K138 (cleared)
$ws.Cells.Item(138, 12).Value = $null  # L138 (cleared)
$ws.Cells.Item(138, 13).Value = $null  # M138 (cleared)
$ws.Cells.Item(138, 14).Value = $null  # N138 (cleared)
$ws.Cells.Item(139, 8).Value = $null  # H139 (cleared)
$ws.Cells.Item(139, 9).Value = $null  # I139 (cleared)
$ws.Cells.Item(139, 10).Value = $null  # J139 (cleared)
$ws.Cells.Item(139, 11).Value = $null  # K139 (cleared)
$ws.Cells.Item(139, 12).Value = $null  # L139 (cleared)
$ws.Cells.Item(140, 8).Value = $null  # H140 (cleared)
$ws.Cells.Item(140, 9).Value = $null  # I140 (cleared)
$ws.Cells.Item(140, 10).Value = $null  # J140 (cleared)
$ws.Cells.Item(140, 11).Value = $null  # K140 (cleared)
$ws.Cells.Item(140, 12).Value = $null  # L140 (cleared)
$ws.Cells.Item(141, 8).Value = $null  # H141 (cleared)
$ws.Cells.Item(141, 9).Value = $null  # I141 (cleared)
$ws.Cells.Item(141, 10).Value = $null  # J141 (cleared)
$ws.Cells.Item(141, 11).Value = $null  # K141 (cleared)
$ws.Cells.Item(141, 12).Value = $null  # L141 (cleared)
$ws.Cells.Item(141, 13).Value = $null  # M141 (cleared)
$ws.Cells.Item(141, 14).Value = $null  # N141 (cleared)
# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(55, 8).Value = 22199.4  # H55
$ws.Cells.Item(55, 10).Value = 22999.5  # J55
$ws.Cells.Item(55, 12).Value = 22999.5  # L55
$ws.Cells.Item(55, 14).Value = -23629.5  # N55
$ws.Cells.Item(74, 8).Value = 1520.5625  # H74
$ws.Cells.Item(74, 9).Value = 1036.8462  # I74
$ws.Cells.Item(74, 10).Value = 3616.6667  # J74
$ws.Cells.Item(74, 11).Value = 1036.8462  # K74
$ws.Cells.Item(74, 12).Value = 3616.6667  # L74
$ws.Cells.Item(74, 13).Value = -162.8462  # M74
$ws.Cells.Item(74, 14).Value = -5364.6667  # N74
$ws.Cells.Item(77, 8).Value = 1520.5625  # H77
$ws.Cells.Item(77, 9).Value = 1036.8462  # I77
$ws.Cells.Item(77, 10).Value = 3616.6667  # J77
$ws.Cells.Item(77, 11).Value = 5184.231  # K77
$ws.Cells.Item(77, 12).Value = 18083.3335  # L77
$ws.Cells.Item(77, 13).Value = -816.2309999999998  # M77
$ws.Cells.Item(77, 14).Value = -26819.3335  # N77
$ws.Cells.Item(121, 8).Value = $null  # H121 (cleared)
$ws.Cells.Item(121, 9).Value = $null  # I121 (cleared)
$ws.Cells.Item(121, 10).Value = $null  # J121 (cleared)
$ws.Cells.Item(121, 11).Value = $null  # K121 (cleared)
$ws.Cells.Item(121, 12).Value = $null  # L121 (cleared)
$ws.Cells.Item(122, 8).Value = $null  # H122 (cleared)
$ws.Cells.Item(122, 9).Value = $null  # I122 (cleared)
$ws.Cells.Item(122, 10).Value = $null  # J122 (cleared)
$ws.Cells.Item(122, 11).Value = $null  # K122 (cleared)
$ws.Cells.Item(122, 12).Value = $null  # L122 (cleared)
$ws.Cells.Item(122, 13).Value = $null  # M122 (cleared)
$ws.Cells.Item(122, 14).Value = $null  # N122 (cleared)
$ws.Cells.Item(123, 8).Value = $null  # H123 (cleared)
$ws.Cells.Item(123, 9).Value = $null  # I123 (cleared)
$ws.Cells.Item(123, 10).Value = $null  # J123 (cleared)
$ws.Cells.Item(123, 11).Value = $null  # K123 (cleared)
$ws.Cells.Item(123, 12).Value = $null  # L123 (cleared)
$ws.Cells.Item(123, 13).Value = $null  # M123 (cleared)
$ws.Cells.Item(123, 14).Value = $null  # N123 (cleared)
$ws.Cells.Item(124, 8).Value = $null  # H124 (cleared)
$ws.Cells.Item(124, 9).Value = $null  # I124 (cleared)
$ws.Cells.Item(124, 10).Value = $null  # J124 (cleared)
$ws.Cells.Item(124, 11).Value = $null  # K124 (cleared)
$ws.Cells.Item(124, 12).Value = $null  # L124 (cleared)
$ws.Cells.Item(124, 14).Value = $null  # N124 (cleared)
$ws.Cells.Item(125, 8).Value = $null  # H125 (cleared)
$ws.Cells.Item(125, 9).Value = $null  # I125 (cleared)
$ws.Cells.Item(125, 10).Value = $null  # J125 (cleared)
$ws.Cells.Item(125, 11).Value = $null  # K125 (cleared)
$ws.Cells.Item(125, 12).Value = $null  # L125 (cleared)
$ws.Cells.Item(125, 14).Value = $null  # N125 (cleared)
$ws.Cells.Item(126, 8).Value = $null  # H126 (cleared)
$ws.Cells.Item(126, 9).Value = $null  # I126 (cleared)
$ws.Cells.Item(126, 10).Value = $null  # J126 (cleared)
$ws.Cells.Item(126, 11).Value = $null  # K126 (cleared)
$ws.Cells.Item(126, 12).Value = $null  # L126 (cleared)
$ws.Cells.Item(127, 8).Value = $null  # H127 (cleared)
$ws.Cells.Item(127, 9).Value = $null  # I127 (cleared)
$ws.Cells.Item(127, 10).Value = $null  # J127 (cleared)
$ws.Cells.Item(127, 11).Value = $null  # K127 (cleared)
$ws.Cells.Item(127, 12).Value = $null  # L127 (cleared)
$ws.Cells.Item(128, 8).Value = $null  # H128 (cleared)
$ws.Cells.Item(128, 9).Value = $null  # I128 (cleared)
$ws.Cells.Item(128, 10).Value = $null  # J128 (cleared)
$ws.Cells.Item(128, 11).Value = $null  # K128 (cleared)
$ws.Cells.Item(128, 12).Value = $null  # L128 (cleared)
$ws.Cells.Item(129, 8).Value = $null  # H129 (cleared)
$ws.Cells.Item(129, 9).Value = $null  # I129 (cleared)
$ws.Cells.Item(129, 10).Value = $null  # J129 (cleared)
$ws.Cells.Item(129, 11).Value = $null  # K129 (cleared)
$ws.Cells.Item(129, 12).Value = $null  # L129 (cleared)
$ws.Cells.Item(130, 8).Value = $null  # H130 (cleared)
$ws.Cells.Item(130, 9).Value = $null  # I130 (cleared)
$ws.Cells.Item(130, 10).Value = $null  # J130 (cleared)
$ws.Cells.Item(130, 11).Value = $null  # K130 (cleared)
$ws.Cells.Item(130, 12).Value = $null  # L130 (cleared)
$ws.Cells.Item(130, 14).Value = $null  # N130 (cleared)
$ws.Cells.Item(131, 8).Value = $null  # H131 (cleared)
$ws.Cells.Item(131, 9).Value = $null  # I131 (cleared)
$ws.Cells.Item(131, 10).Value = $null  # J131 (cleared)
$ws.Cells.Item(131, 11).Value = $null  # K131 (cleared)
$ws.Cells.Item(131, 12).Value = $null  # L131 (cleared)
$ws.Cells.Item(132, 8).Value = $null  # H132 (cleared)
$ws.Cells.Item(132, 9).Value = $null  # I132 (cleared)
$ws.Cells.Item(132, 10).Value = $null  # J132 (cleared)
$ws.Cells.Item(132, 11).Value = $null  # K132 (cleared)
$ws.Cells.Item(132, 12).Value = $null  # L132 (cleared)
$ws.Cells.Item(132, 13).Value = $null  # M132 (cleared)
$ws.Cells.Item(132, 14).Value = $null  # N132 (cleared)
$ws.Cells.Item(133, 8).Value = $null  # H133 (cleared)
$ws.Cells.Item(133, 9).Value = $null  # I133 (cleared)
$ws.Cells.Item(133, 10).Value = $null  # J133 (cleared)
$ws.Cells.Item(133, 11).Value = $null  # K133 (cleared)
$ws.Cells.Item(133, 12).Value = $null  # L133 (cleared)
$ws.Cells.Item(134, 8).Value = $null  # H134 (cleared)
$ws.Cells.Item(134, 9).Value = $null  # I134 (cleared)
$ws.Cells.Item(134, 10).Value = $null  # J134 (cleared)
$ws.Cells.Item(134, 11).Value = $null  # K134 (cleared)
$ws.Cells.Item(134, 12).Value = $null  # L134 (cleared)
$ws.Cells.Item(134, 14).Value = $null  # N134 (cleared)
$ws.Cells.Item(135, 8).Value = $null  # H135 (cleared)
$ws.Cells.Item(135, 9).Value = $null  # I135 (cleared)
$ws.Cells.Item(135, 10).Value = $null  # J135 (cleared)
$ws.Cells.Item(135, 11).Value = $null  # K135 (cleared)
$ws.Cells.Item(135, 12).Value = $null  # L135 (cleared)
$ws.Cells.Item(135, 14).Value = $null  # N135 (cleared)
$ws.Cells.Item(137, 8).Value = $null  # H137 (cleared)
$ws.Cells.Item(137, 9).Value = $null  # I137 (cleared)
$ws.Cells.Item(137, 10).Value = $null  # J137 (cleared)
$ws.Cells.Item(137, 11).Value = $null  # K137 (cleared)
$ws.Cells.Item(137, 12).Value = $null  # L137 (cleared)
$ws.Cells.Item(138, 8).Value = $null  # H138 (cleared)
$ws.Cells.Item(138, 9).Value = $null  # I138 (cleared)
$ws.Cells.Item(138, 10).Value = $null  # J138 (cleared)
$ws.Cells.Item(138, 11).Value = $null  # K138 (cleared)
$ws.Cells.Item(138, 12).Value = $null  # L138 (cleared)
$ws.Cells.Item(139, 8).Value = $null  # H139 (cleared)
$ws.Cells.Item(139, 9).Value = $null  # I139 (cleared)
$ws.Cells.Item(139, 10).Value = $null  # J139 (cleared)
$ws.Cells.Item(139, 11).Value = $null  # K139 (cleared)
$ws.Cells.Item(139, 12).Value = $null  # L139 (cleared)
$ws.Cells.Item(139, 13).Value = $null  # M139 (cleared)
$ws.Cells.Item(140, 8).Value = $null  # H140 (cleared)
$ws.Cells.Item(140, 9).Value = $null  # I140 (cleared)
$ws.Cells.Item(140, 10).Value = $null  # J140 (cleared)
$ws.Cells.Item(140, 11).Value = $null  # K140 (cleared)
$ws.Cells.Item(140, 12).Value = $null  # L140 (cleared)
$ws.Cells.Item(141, 8).Value = $null  # H141 (cleared)
$ws.Cells.Item(141, 9).Value = $null  # I141 (cleared)
$ws.Cells.Item(141, 10).Value = $null  # J141 (cleared)
$ws.Cells.Item(141, 11).Value = $null  # K141 (cleared)
$ws.Cells.Item(141, 12).Value = $null  # L141 (cleared)
# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2098.75  # H99
$ws.Cells.Item(99, 9).Value = 1845.2  # I99
$ws.Cells.Item(99, 11).Value = 1845.2  # K99
$ws.Cells.Item(99, 13).Value = -347.2  # M99
# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 478.69696  # H7
$ws.Cells.Item(7, 9).Value = 438.6842  # I7
$ws.Cells.Item(7, 10).Value = 533  # J7
$ws.Cells.Item(7, 11).Value = 438.6842  # K7
$ws.Cells.Item(7, 12).Value = 533  # L7
$ws.Cells.Item(7, 13).Value = -325.6842  # M7
$ws.Cells.Item(7, 14).Value = -759  # N7
$ws.Cells.Item(17, 8).Value = 2001  # H17
$ws.Cells.Item(17, 9).Value = 0  # I17
$ws.Cells.Item(17, 11).Value = 0  # K17
$ws.Cells.Item(17, 13).Value = $null  # M17 (cleared)
$ws.Cells.Item(19, 8).Value = 753.8889  # H19
$ws.Cells.Item(19, 9).Value = 966.9231  # I19
$ws.Cells.Item(19, 11).Value = 966.9231  # K19
$ws.Cells.Item(19, 13).Value = -796.9231  # M19
$ws.Cells.Item(24, 8).Value = 753.8889  # H24
$ws.Cells.Item(24, 9).Value = 966.9231  # I24
$ws.Cells.Item(24, 11).Value = 966.9231  # K24
$ws.Cells.Item(24, 13).Value = -796.9231  # M24
$ws.Cells.Item(31, 8).Value = 3668.1738  # H31
$ws.Cells.Item(31, 9).Value = 4096.375  # I31
$ws.Cells.Item(31, 10).Value = 3439.8  # J31
$ws.Cells.Item(31, 11).Value = 4096.375  # K31
$ws.Cells.Item(31, 12).Value = 3439.8  # L31
$ws.Cells.Item(31, 13).Value = -3801.375  # M31
$ws.Cells.Item(31, 14).Value = -4029.8  # N31
$ws.Cells.Item(34, 8).Value = 3668.1738  # H34
$ws.Cells.Item(34, 9).Value = 4096.375  # I34
$ws.Cells.Item(34, 10).Value = 3439.8  # J34
$ws.Cells.Item(34, 11).Value = 4096.375  # K34
$ws.Cells.Item(34, 12).Value = 3439.8  # L34
$ws.Cells.Item(34, 13).Value = -3894.375  # M34
$ws.Cells.Item(34, 14).Value = -3843.8  # N34
$ws.Cells.Item(41, 8).Value = 13169.923  # H41
$ws.Cells.Item(41, 9).Value = 12253.6  # I41
$ws.Cells.Item(41, 10).Value = 13742.625  # J41
$ws.Cells.Item(41, 11).Value = 12253.6  # K41
$ws.Cells.Item(41, 12).Value = 13742.625  # L41
$ws.Cells.Item(41, 13).Value = -11825.6  # M41
$ws.Cells.Item(41, 14).Value = -14598.625  # N41
$ws.Cells.Item(50, 8).Value = 18199.8  # H50
$ws.Cells.Item(50, 10).Value = 18999.75  # J50
$ws.Cells.Item(50, 12).Value = 18999.75  # L50
$ws.Cells.Item(50, 14).Value = -20249.75  # N50
$ws.Cells.Item(59, 8).Value = 28375  # H59
$ws.Cells.Item(59, 9).Value = 23500  # I59
$ws.Cells.Item(59, 11).Value = 23500  # K59
$ws.Cells.Item(59, 13).Value = -22355  # M59
$ws.Cells.Item(60, 8).Value = 21999.75  # H60
$ws.Cells.Item(60, 9).Value = 22666.334  # I60
$ws.Cells.Item(60, 10).Value = 20000  # J60
$ws.Cells.Item(60, 11).Value = 22666.334  # K60
$ws.Cells.Item(60, 12).Value = 20000  # L60
$ws.Cells.Item(60, 13).Value = -22155.334  # M60
$ws.Cells.Item(60, 14).Value = -21022  # N60
$ws.Cells.Item(129, 8).Value = 0  # H129
$ws.Cells.Item(129, 9).Value = 0  # I129
$ws.Cells.Item(129, 10).Value = 0  # J129
$ws.Cells.Item(129, 11).Value = 0  # K129
$ws.Cells.Item(129, 12).Value = 0  # L129
$ws.Cells.Item(130, 8).Value = 0  # H130
$ws.Cells.Item(130, 9).Value = 0  # I130
$ws.Cells.Item(130, 10).Value = 0  # J130
$ws.Cells.Item(130, 11).Value = 0  # K130
$ws.Cells.Item(130, 12).Value = 0  # L130
$ws.Cells.Item(131, 8).Value = 50000  # H131
$ws.Cells.Item(131, 9).Value = 0  # I131
$ws.Cells.Item(131, 10).Value = 50000  # J131
$ws.Cells.Item(131, 11).Value = 0  # K131
$ws.Cells.Item(131, 12).Value = 50000  # L131
$ws.Cells.Item(131, 14).Value = -60080  # N131
$ws.Cells.Item(132, 8).Value = 2673.625  # H132
$ws.Cells.Item(132, 9).Value = 2180.75  # I132
$ws.Cells.Item(132, 10).Value = 4152.25  # J132
$ws.Cells.Item(132, 11).Value = 6542.25  # K132
$ws.Cells.Item(132, 12).Value = 12456.75  # L132
$ws.Cells.Item(132, 13).Value = -4012.25  # M132
$ws.Cells.Item(132, 14).Value = -17516.75  # N132
$ws.Cells.Item(133, 8).Value = 99000  # H133
$ws.Cells.Item(133, 9).Value = 0  # I133
$ws.Cells.Item(133, 10).Value = 99000  # J133
$ws.Cells.Item(133, 11).Value = 0  # K133
$ws.Cells.Item(133, 12).Value = 99000  # L133
$ws.Cells.Item(133, 14).Value = -104060  # N133
$ws.Cells.Item(134, 8).Value = 3657.4546  # H134
$ws.Cells.Item(134, 9).Value = 3529.375  # I134
$ws.Cells.Item(134, 10).Value = 3999  # J134
$ws.Cells.Item(134, 11).Value = 10588.125  # K134
$ws.Cells.Item(134, 12).Value = 11997  # L134
$ws.Cells.Item(134, 13).Value = -8053.125  # M134
$ws.Cells.Item(134, 14).Value = -17067  # N134
$ws.Cells.Item(135, 8).Value = 86999.5  # H135
$ws.Cells.Item(135, 9).Value = 0  # I135
$ws.Cells.Item(135, 10).Value = 86999.5  # J135
$ws.Cells.Item(135, 11).Value = 0  # K135
$ws.Cells.Item(135, 12).Value = 86999.5  # L135
$ws.Cells.Item(135, 14).Value = -97139.5  # N135
$ws.Cells.Item(137, 8).Value = 55000  # H137
$ws.Cells.Item(137, 9).Value = 0  # I137
$ws.Cells.Item(137, 10).Value = 55000  # J137
$ws.Cells.Item(137, 11).Value = 0  # K137
$ws.Cells.Item(137, 12).Value = 55000  # L137
$ws.Cells.Item(137, 14).Value = -65200  # N137
$ws.Cells.Item(138, 8).Value = 100000  # H138
$ws.Cells.Item(138, 9).Value = 0  # I138
$ws.Cells.Item(138, 10).Value = 100000  # J138
$ws.Cells.Item(138, 11).Value = 0  # K138
$ws.Cells.Item(138, 12).Value = 100000  # L138
$ws.Cells.Item(138, 14).Value = -110280  # N138
$ws.Cells.Item(139, 8).Value = 69999.5  # H139
$ws.Cells.Item(139, 9).Value = 0  # I139
$ws.Cells.Item(139, 10).Value = 69999.5  # J139
$ws.Cells.Item(139, 11).Value = 0  # K139
$ws.Cells.Item(139, 12).Value = 69999.5  # L139
$ws.Cells.Item(139, 14).Value = -80279.5  # N139
$ws.Cells.Item(140, 8).Value = 60000  # H140
$ws.Cells.Item(140, 9).Value = 0  # I140
$ws.Cells.Item(140, 10).Value = 60000  # J140
$ws.Cells.Item(140, 11).Value = 0  # K140
$ws.Cells.Item(140, 12).Value = 60000  # L140
$ws.Cells.Item(140, 14).Value = -70360  # N140
$ws.Cells.Item(141, 8).Value = 82220.5  # H141
$ws.Cells.Item(141, 9).Value = 0  # I141
$ws.Cells.Item(141, 10).Value = 82220.5  # J141
$ws.Cells.Item(141, 11).Value = 0  # K141
$ws.Cells.Item(141, 12).Value = 82220.5  # L141
$ws.Cells.Item(141, 14).Value = -92580.5  # N141
# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 631.3333  # H5
$ws.Cells.Item(5, 9).Value = 0  # I5
$ws.Cells.Item(5, 10).Value = 631.3333  # J5
$ws.Cells.Item(5, 11).Value = 0  # K5
$ws.Cells.Item(5, 12).Value = $null  # L5 (cleared)
$ws.Cells.Item(5, 13).Value = 1893.9999  # M5
$ws.Cells.Item(5, 14).Value = -2117.9999  # N5
$ws.Cells.Item(68, 8).Value = 4291.773  # H68
$ws.Cells.Item(68, 9).Value = 2450  # I68
$ws.Cells.Item(68, 10).Value = 4701.0557  # J68
$ws.Cells.Item(68, 11).Value = 7350  # K68
$ws.Cells.Item(68, 12).Value = 14103.1671  # L68
$ws.Cells.Item(68, 13).Value = -6539  # M68
$ws.Cells.Item(68, 14).Value = -15725.1671  # N68
$ws.Cells.Item(71, 8).Value = 4291.773  # H71
$ws.Cells.Item(71, 9).Value = 2450  # I71
$ws.Cells.Item(71, 10).Value = 4701.0557  # J71
$ws.Cells.Item(71, 11).Value = 22050  # K71
$ws.Cells.Item(71, 12).Value = 42309.5013  # L71
$ws.Cells.Item(71, 13).Value = -17994  # M71
$ws.Cells.Item(71, 14).Value = -50421.5013  # N71
$ws.Cells.Item(135, 8).Value = 631.3333  # H135
$ws.Cells.Item(135, 9).Value = 0  # I135
$ws.Cells.Item(135, 10).Value = 631.3333  # J135
$ws.Cells.Item(135, 11).Value = 0  # K135
$ws.Cells.Item(135, 12).Value = $null  # L135 (cleared)
$ws.Cells.Item(135, 13).Value = 5681.9997  # M135
$ws.Cells.Item(135, 14).Value = -10751.9997  # N135
# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 4569.5  # H126
$ws.Cells.Item(126, 9).Value = 4386.4546  # I126
$ws.Cells.Item(126, 10).Value = 4857.143  # J126
$ws.Cells.Item(126, 11).Value = 13159.3638  # K126
$ws.Cells.Item(126, 12).Value = 14571.429  # L126
$ws.Cells.Item(126, 13).Value = -10689.3638  # M126
$ws.Cells.Item(126, 14).Value = -19511.429  # N126
# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 20000  # H12
$ws.Cells.Item(12, 9).Value = 20000  # I12
$ws.Cells.Item(12, 10).Value = 0  # J12
$ws.Cells.Item(12, 11).Value = 20000  # K12
$ws.Cells.Item(12, 12).Value = 0  # L12
$ws.Cells.Item(12, 13).Value = $null  # M12 (cleared)
$ws.Cells.Item(12, 14).Value = -19858  # N12
$ws.Cells.Item(119, 8).Value = $null  # H119 (cleared)
$ws.Cells.Item(119, 9).Value = $null  # I119 (cleared)
$ws.Cells.Item(119, 10).Value = $null  # J119 (cleared)
$ws.Cells.Item(119, 11).Value = $null  # K119 (cleared)
$ws.Cells.Item(119, 12).Value = $null  # L119 (cleared)
$ws.Cells.Item(119, 13).Value = $null  # M119 (cleared)
$ws.Cells.Item(119, 14).Value = $null  # N119 (cleared)
$ws.Cells.Item(120, 8).Value = $null  # H120 (cleared)
$ws.Cells.Item(120, 9).Value = $null  # I120 (cleared)
$ws.Cells.Item(120, 10).Value = $null  # J120 (cleared)
$ws.Cells.Item(120, 11).Value = $null  # K120 (cleared)
$ws.Cells.Item(120, 12).Value = $null  # L120 (cleared)
$ws.Cells.Item(121, 8).Value = $null  # H121 (cleared)
$ws.Cells.Item(121, 9).Value = $null  # I121 (cleared)
$ws.Cells.Item(121, 10).Value = $null  # J121 (cleared)
$ws.Cells.Item(121, 11).Value = $null  # K121 (cleared)
$ws.Cells.Item(121, 12).Value = $null  # L121 (cleared)
$ws.Cells.Item(122, 8).Value = $null  # H122 (cleared)
$ws.Cells.Item(122, 9).Value = $null  # I122 (cleared)
$ws.Cells.Item(122, 10).Value = $null  # J122 (cleared)
$ws.Cells.Item(122, 11).Value = $null  # K122 (cleared)
$ws.Cells.Item(122, 12).Value = $null  # L122 (cleared)
$ws.Cells.Item(122, 13).Value = $null  # M122 (cleared)
$ws.Cells.Item(122, 14).Value = $null  # N122 (cleared)
$ws.Cells.Item(123, 8).Value = $null  # H123 (cleared)
$ws.Cells.Item(123, 9).Value = $null  # I123 (cleared)
$ws.Cells.Item(123, 10).Value = $null  # J123 (cleared)
$ws.Cells.Item(123, 11).Value = $null  # K123 (cleared)
$ws.Cells.Item(123, 12).Value = $null  # L123 (cleared)
$ws.Cells.Item(124, 8).Value = $null  # H124 (cleared)
$ws.Cells.Item(124, 9).Value = $null  # I124 (cleared)
$ws.Cells.Item(124, 10).Value = $null  # J124 (cleared)
$ws.Cells.Item(124, 11).Value = $null  # K124 (cleared)
$ws.Cells.Item(124, 12).Value = $null  # L124 (cleared)
$ws.Cells.Item(125, 8).Value = $null  # H125 (cleared)
$ws.Cells.Item(125, 9).Value = $null  # I125 (cleared)
$ws.Cells.Item(125, 10).Value = $null  # J125 (cleared)
$ws.Cells.Item(125, 11).Value = $null  # K125 (cleared)
$ws.Cells.Item(125, 12).Value = $null  # L125 (cleared)
$ws.Cells.Item(125, 14).Value = $null  # N125 (cleared)
$ws.Cells.Item(126, 8).Value = $null  # H126 (cleared)
$ws.Cells.Item(126, 9).Value = $null  # I126 (cleared)
$ws.Cells.Item(126, 10).Value = $null  # J126 (cleared)
$ws.Cells.Item(126, 11).Value = $null  # K126 (cleared)
$ws.Cells.Item(126, 12).Value = $null  # L126 (cleared)
$ws.Cells.Item(126, 13).Value = $null  # M126 (cleared)
$ws.Cells.Item(126, 14).Value = $null  # N126 (cleared)
$ws.Cells.Item(127, 8).Value = $null  # H127 (cleared)
$ws.Cells.Item(127, 9).Value = $null  # I127 (cleared)
$ws.Cells.Item(127, 10).Value = $null  # J127 (cleared)
$ws.Cells.Item(127, 11).Value = $null  # K127 (cleared)
$ws.Cells.Item(127, 12).Value = $null  # L127 (cleared)
$ws.Cells.Item(127, 14).Value = $null  # N127 (cleared)
$ws.Cells.Item(128, 8).Value = $null  # H128 (cleared)
$ws.Cells.Item(128, 9).Value = $null  # I128 (cleared)
$ws.Cells.Item(128, 10).Value = $null  # J128 (cleared)
$ws.Cells.Item(128, 11).Value = $null  # K128 (cleared)
$ws.Cells.Item(128, 12).Value = $null  # L128 (cleared)
$ws.Cells.Item(129, 8).Value = $null  # H129 (cleared)
$ws.Cells.Item(129, 9).Value = $null  # I129 (cleared)
$ws.Cells.Item(129, 10).Value = $null  # J129 (cleared)
$ws.Cells.Item(129, 11).Value = $null  # K129 (cleared)
$ws.Cells.Item(129, 12).Value = $null  # L129 (cleared)
$ws.Cells.Item(130, 8).Value = $null  # H130 (cleared)
$ws.Cells.Item(130, 9).Value = $null  # I130 (cleared)
$ws.Cells.Item(130, 10).Value = $null  # J130 (cleared)
$ws.Cells.Item(130, 11).Value = $null  # K130 (cleared)
$ws.Cells.Item(130, 12).Value = $null  # L130 (cleared)
$ws.Cells.Item(131, 8).Value = $null  # H131 (cleared)
$ws.Cells.Item(131, 9).Value = $null  # I131 (cleared)
$ws.Cells.Item(131, 10).Value = $null  # J131 (cleared)
$ws.Cells.Item(131, 11).Value = $null  # K131 (cleared)
$ws.Cells.Item(131, 12).Value = $null  # L131 (cleared)
$ws.Cells.Item(131, 14).Value = $null  # N131 (cleared)
$ws.Cells.Item(132, 8).Value = $null  # H132 (cleared)
$ws.Cells.Item(132, 9).Value = $null  # I132 (cleared)
$ws.Cells.Item(132, 10).Value = $null  # J132 (cleared)
$ws.Cells.Item(132, 11).Value = $null  # K132 (cleared)
$ws.Cells.Item(132, 12).Value = $null  # L132 (cleared)
$ws.Cells.Item(132, 13).Value = $null  # M132 (cleared)
$ws.Cells.Item(132, 14).Value = $null  # N132 (cleared)
$ws.Cells.Item(133, 8).Value = $null  # H133 (cleared)
$ws.Cells.Item(133, 9).Value = $null  # I133 (cleared)
$ws.Cells.Item(133, 10).Value = $null  # J133 (cleared)
$ws.Cells.Item(133, 11).Value = $null  # K133 (cleared)
$ws.Cells.Item(133, 12).Value = $null  # L133 (cleared)
$ws.Cells.Item(133, 14).Value = $null  # N133 (cleared)
$ws.Cells.Item(135, 8).Value = $null  # H135 (cleared)
$ws.Cells.Item(135, 9).Value = $null  # I135 (cleared)
$ws.Cells.Item(135, 10).Value = $null  # J135 (cleared)
$ws.Cells.Item(135, 11).Value = $null  # K135 (cleared)
$ws.Cells.Item(135, 12).Value = $null  # L135 (cleared)
$ws.Cells.Item(136, 8).Value = $null  # H136 (cleared)
$ws.Cells.Item(136, 9).Value = $null  # I136 (cleared)
$ws.Cells.Item(136, 10).Value = $null  # J136 (cleared)
$ws.Cells.Item(136, 11).Value = $null  # K136 (cleared)
$ws.Cells.Item(136, 12).Value = $null  # L136 (cleared)
$ws.Cells.Item(136, 13).Value = $null  # M136 (cleared)
$ws.Cells.Item(137, 8).Value = $null  # H137 (cleared)
$ws.Cells.Item(137, 9).Value = $null  # I137 (cleared)
$ws.Cells.Item(137, 10).Value = $null  # J137 (cleared)
$ws.Cells.Item(137, 11).Value = $null  # K137 (cleared)
$ws.Cells.Item(137, 12).Value = $null  # L137 (cleared)
$ws.Cells.Item(137, 14).Value = $null  # N137 (cleared)
$ws.Cells.Item(138, 8).Value = $null  # H138 (cleared)
$ws.Cells.Item(138, 9).Value = $null  # I138 (cleared)
$ws.Cells.Item(138, 10).Value = $null  # J138 (cleared)
$ws.Cells.Item(138, 11).Value = $null  # K138 (cleared)
$ws.Cells.Item(138, 12).Value = $null  # L138 (cleared)
$ws.Cells.Item(139, 8).Value = $null  # H139 (cleared)
$ws.Cells.Item(139, 9).Value = $null  # I139 (cleared)
$ws.Cells.Item(139, 10).Value = $null  # J139 (cleared)
$ws.Cells.Item(139, 11).Value = $null  # K139 (cleared)
$ws.Cells.Item(139, 12).Value = $null  # L139 (cleared)
$ws.Cells.Item(140, 8).Value = $null  # H140 (cleared)
$ws.Cells.Item(140, 9).Value = $null  # I140 (cleared)
$ws.Cells.Item(140, 10).Value = $null  # J140 (cleared)
$ws.Cells.Item(140, 11).Value = $null  # K140 (cleared)
$ws.Cells.Item(140, 12).Value = $null  # L140 (cleared)
$ws.Cells.Item(141, 8).Value = $null  # H141 (cleared)
$ws.Cells.Item(141, 9).Value = $null  # I141 (cleared)
$ws.Cells.Item(141, 10).Value = $null  # J141 (cleared)
$ws.Cells.Item(141, 11).Value = $null  # K141 (cleared)
$ws.Cells.Item(141, 12).Value = $null  # L141 (cleared)
